$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.664.24'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = '1.739.35'
$ws.Range('E3').Value = '  -5.79%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -10.04%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4898'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -7.92%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.34'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -8.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2551'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -18.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06021'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -12.72%  '
$ws.Range('D11').Value = '1.743.18'
$ws.Range('E11').Value = '  -5.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.06815'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -12.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.79'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -20.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.417'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -12.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.12'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -14.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.5677'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -25.56%  '
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('D19').Value = '25.711.03'
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -19.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000006544'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -17.60%  '
$ws.Range('D22').Value = '1.965.07'
$ws.Range('E22').Value = '  -5.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.017'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -13.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.894'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -15.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.014'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -16.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '136.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.68%  '
$ws.Range('E27').Value = '  -12.47%  '
$ws.Range('E28').Value = '  -17.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -14.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '101.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -9.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.726'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -13.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07949'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -9.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.367'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -17.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04376'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.001'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E36').Value = '  -11.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9724'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -14.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5934'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -19.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.658'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -14.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.18'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.3715'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -22.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05217'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -10.28%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1069'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -14.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '29.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -14.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '51.87'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -14.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.772'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -24.51%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.001'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.04%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01509'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -12.68%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.876'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -19.20%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.138'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -12.90%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7444'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -17.87%  '
